$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 219 (pushes existing rows 219..319 down to 220..320)
$ws.Rows.Item(219).Insert()

# Populate the new row with the record's values
$ws.Cells.Item(219, 1).Value = 10
$ws.Cells.Item(219, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(219, 3).Value = "La Araucanía"
$ws.Cells.Item(219, 4).Value = 44553
$ws.Cells.Item(219, 5).Value = 9
$ws.Cells.Item(219, 6).Value = 100114014
$ws.Cells.Item(219, 7).Value = "Betarraga"
$ws.Cells.Item(219, 8).Value = "Sin especificar"
$ws.Cells.Item(219, 9).Value = "Primera"
$ws.Cells.Item(219, 10).Value = 125
$ws.Cells.Item(219, 11).Value = 700
$ws.Cells.Item(219, 12).Value = 700
$ws.Cells.Item(219, 13).Value = 700
$ws.Cells.Item(219, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(219, 15).Value = "Región del Maule"
$ws.Cells.Item(219, 16).Value = 140
$ws.Cells.Item(219, 17).Value = 5
$ws.Cells.Item(219, 18).Value = "Hortaliza"
